# Apply the Alvearie FHIR IG deploy-refresh edit:
#  - bump Version / Date on the Metadata sheet
#  - replace the duplicated "Contact" rows with Publisher/Jurisdiction info
#  - update Short/Definition for the root Extension row on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Rows 10 and 11 both currently hold "Contact" / "No display for ContactDetail".
# Remove the duplicate row 11 entirely, then turn the remaining pair into the
# Publisher value (row 9) and a new Jurisdiction row (row 10).
$meta.Rows.Item(11).Delete()

$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Elements sheet: update the Short / Definition text for the root Extension row
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Episode Primary Physician"
$elements.Range("L2").Value = "Primary physician associated with the episode of care"
